# Swap the "System, " / "<email>" ordering in the "Recorded By" column (G)
# so that existing cells reading "System, dnasr281@gmail.com" become
# "dnasr281@gmail.com, System". Cells containing only one of the two
# values (e.g. just "System" or just the email) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
}
